$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 965
$ws.Range("I18").Value = 965
$ws.Range("K18").Value = 965
$ws.Range("M18").Value = -681

$ws.Range("H43").Value = 1591.2632
$ws.Range("J43").Value = 1367.9231
$ws.Range("L43").Value = 1367.9231
$ws.Range("N43").Value = -1505.9231

$ws.Range("H62").Value = 4947.091
$ws.Range("I62").Value = 5936.8335
$ws.Range("K62").Value = 5936.8335
$ws.Range("M62").Value = -5312.8335

$ws.Range("H65").Value = 4947.091
$ws.Range("I65").Value = 5936.8335
$ws.Range("K65").Value = 29684.1675
$ws.Range("M65").Value = -26564.1675

$ws.Range("H70").Value = 3259.1035
$ws.Range("I70").Value = 1697.7142
$ws.Range("J70").Value = 3755.9092
$ws.Range("K70").Value = 5093.142599999999
$ws.Range("L70").Value = 11267.7276
$ws.Range("M70").Value = -4823.142599999999
$ws.Range("N70").Value = -11807.7276

$ws.Range("H73").Value = 3259.1035
$ws.Range("I73").Value = 1697.7142
$ws.Range("J73").Value = 3755.9092
$ws.Range("K73").Value = 5093.142599999999
$ws.Range("L73").Value = 11267.7276
$ws.Range("M73").Value = -4157.142599999999
$ws.Range("N73").Value = -13139.7276

$ws.Range("H98").Value = 4334.905
$ws.Range("I98").Value = 2449.1052
$ws.Range("K98").Value = 2449.1052
$ws.Range("M98").Value = -951.1052

$ws.Range("H113").Value = 17422.643
$ws.Range("I113").Value = 37334.5
$ws.Range("J113").Value = 2488.75
$ws.Range("K113").Value = 37334.5
$ws.Range("L113").Value = 2488.75
$ws.Range("M113").Value = -34080.5
$ws.Range("N113").Value = -8996.75

$ws.Range("H122").Value = 4334.905
$ws.Range("I122").Value = 2449.1052
$ws.Range("K122").Value = 7347.3156
$ws.Range("M122").Value = -4897.3156

$ws.Range("H137").Value = 1872.625
$ws.Range("I137").Value = 996.2
$ws.Range("J137").Value = 3333.3333
$ws.Range("K137").Value = 2988.6
$ws.Range("L137").Value = 9999.999899999999
$ws.Range("M137").Value = -438.6000000000004
$ws.Range("N137").Value = -15099.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2611
$ws.Range("I2").Value = 1873.9166
$ws.Range("J2").Value = 3874.5715
$ws.Range("K2").Value = 1873.9166
$ws.Range("L2").Value = 3874.5715
$ws.Range("M2").Value = -1760.9166
$ws.Range("N2").Value = -4100.5715

$ws.Range("H110").Value = 1472.9131
$ws.Range("I110").Value = 1488.8572
$ws.Range("J110").Value = 1305.5
$ws.Range("K110").Value = 1488.8572
$ws.Range("L110").Value = 1305.5
$ws.Range("M110").Value = 556.1428000000001
$ws.Range("N110").Value = -5395.5

$ws.Range("H116").Value = 2611
$ws.Range("I116").Value = 1873.9166
$ws.Range("J116").Value = 3874.5715
$ws.Range("K116").Value = 1873.9166
$ws.Range("L116").Value = 3874.5715
$ws.Range("M116").Value = 420.0834
$ws.Range("N116").Value = -8462.5715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2611
$ws.Range("I3").Value = 1873.9166
$ws.Range("J3").Value = 3874.5715
$ws.Range("K3").Value = 1873.9166
$ws.Range("L3").Value = 3874.5715
$ws.Range("M3").Value = -1759.9166
$ws.Range("N3").Value = -4102.5715

$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").ClearContents()
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = 0

$ws.Range("H105").Value = 3261.0386
$ws.Range("I105").Value = 3736.6875
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 3736.6875
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = -1989.6875
$ws.Range("N105").Value = -5994

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2591.2
$ws.Range("I16").Value = 5203.6
$ws.Range("K16").Value = 5203.6
$ws.Range("M16").Value = -4916.6

$ws.Range("H31").Value = 5708.75
$ws.Range("I31").Value = 2712.6667
$ws.Range("K31").Value = 2712.6667
$ws.Range("M31").Value = -2417.6667

$ws.Range("H34").Value = 5708.75
$ws.Range("I34").Value = 2712.6667
$ws.Range("K34").Value = 2712.6667
$ws.Range("M34").Value = -2510.6667

$ws.Range("H41").Value = 39163
$ws.Range("I41").Value = 31327
$ws.Range("K41").Value = 31327
$ws.Range("M41").Value = -30899

$ws.Range("H50").Value = 74999.5
$ws.Range("J50").Value = 74999.5
$ws.Range("L50").Value = 74999.5
$ws.Range("N50").Value = -76249.5

$ws.Range("H51").Value = 46598.4
$ws.Range("J51").Value = 56998.332
$ws.Range("L51").Value = 56998.332
$ws.Range("N51").Value = -58470.332

$ws.Range("H61").Value = 46598.4
$ws.Range("J61").Value = 56998.332
$ws.Range("L61").Value = 56998.332
$ws.Range("N61").Value = -57694.332

$ws.Range("H99").Value = 4997
$ws.Range("I99").Value = 4997
$ws.Range("K99").Value = 4997
$ws.Range("M99").Value = -3499

$ws.Range("H113").Value = 2591.2
$ws.Range("I113").Value = 5203.6
$ws.Range("K113").Value = 5203.6
$ws.Range("M113").Value = -3033.6

$ws.Range("H126").Value = 4997
$ws.Range("I126").Value = 4997
$ws.Range("K126").Value = 14991
$ws.Range("M126").Value = -12521

$ws.Range("H134").Value = 2406.4
$ws.Range("I134").Value = 2266.5862
$ws.Range("J134").Value = 3082.1667
$ws.Range("K134").Value = 6799.758600000001
$ws.Range("L134").Value = 9246.500100000001
$ws.Range("M134").Value = -4264.758600000001
$ws.Range("N134").Value = -14316.5001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 2190.889
$ws.Range("I17").Value = 2559.6428
$ws.Range("J17").Value = 900.25
$ws.Range("K17").Value = 7678.928400000001
$ws.Range("L17").Value = 2700.75
$ws.Range("M17").Value = -7509.928400000001
$ws.Range("N17").Value = -3038.75

$ws.Range("H32").Value = 300960.4
$ws.Range("I32").Value = 334267.34
$ws.Range("K32").Value = 1002802.02
$ws.Range("M32").Value = -1002519.02

$ws.Range("H46").Value = 10094388
$ws.Range("I46").Value = 22222412
$ws.Range("K46").Value = 66667236
$ws.Range("M46").Value = -66667145

$ws.Range("H139").Value = 5554.4
$ws.Range("I139").Value = 3291.25
$ws.Range("K139").Value = 9873.75
$ws.Range("M139").Value = -4733.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3648.625
$ws.Range("I80").Value = 3698.3333
$ws.Range("J80").Value = 3499.5
$ws.Range("K80").Value = 3698.3333
$ws.Range("L80").Value = 3499.5
$ws.Range("M80").Value = -2700.3333
$ws.Range("N80").Value = -5495.5

$ws.Range("H83").Value = 3648.625
$ws.Range("I83").Value = 3698.3333
$ws.Range("J83").Value = 3499.5
$ws.Range("K83").Value = 18491.6665
$ws.Range("L83").Value = 17497.5
$ws.Range("M83").Value = -13499.6665
$ws.Range("N83").Value = -27481.5

$ws.Range("H102").Value = 5750.095
$ws.Range("I102").Value = 4455.2
$ws.Range("J102").Value = 6927.273
$ws.Range("K102").Value = 4455.2
$ws.Range("L102").Value = 6927.273
$ws.Range("M102").Value = -2833.2
$ws.Range("N102").Value = -10171.273

$ws.Range("H122").Value = 6665.6665
$ws.Range("I122").Value = 4998.5
$ws.Range("K122").Value = 14995.5
$ws.Range("M122").Value = -12545.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8579.833000000001
$ws.Range("I40").Value = 7595.8
$ws.Range("J40").Value = 13500
$ws.Range("K40").Value = 7595.8
$ws.Range("L40").Value = 13500
$ws.Range("M40").Value = -7459.8
$ws.Range("N40").Value = -13772

$ws.Range("H68").Value = 2385.484
$ws.Range("I68").Value = 1550.6666
$ws.Range("J68").Value = 3541.3845
$ws.Range("K68").Value = 1550.6666
$ws.Range("L68").Value = 3541.3845
$ws.Range("M68").Value = -801.6666
$ws.Range("N68").Value = -5039.3845

$ws.Range("H71").Value = 2385.484
$ws.Range("I71").Value = 1550.6666
$ws.Range("J71").Value = 3541.3845
$ws.Range("K71").Value = 7753.333000000001
$ws.Range("L71").Value = 17706.9225
$ws.Range("M71").Value = -4009.333000000001
$ws.Range("N71").Value = -25194.9225

$ws.Range("H93").Value = 12764
$ws.Range("I93").Value = 26500
$ws.Range("J93").Value = 8185.3335
$ws.Range("K93").Value = 26500
$ws.Range("L93").Value = 8185.3335
$ws.Range("M93").Value = -25252
$ws.Range("N93").Value = -10681.3335

$ws.Range("H122").Value = 6713.2856
$ws.Range("I122").Value = 7197.8
$ws.Range("K122").Value = 21593.4
$ws.Range("M122").Value = -19143.4

$ws.Range("H132").Value = 7698.9185
$ws.Range("I132").Value = 7058.75
$ws.Range("J132").Value = 9471.691999999999
$ws.Range("K132").Value = 21176.25
$ws.Range("L132").Value = 28415.076
$ws.Range("M132").Value = -18646.25
$ws.Range("N132").Value = -33475.076

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 13516822
$ws.Range("I122").Value = 16669951
$ws.Range("K122").Value = 50009853
$ws.Range("M122").Value = -50007403

$ws.Range("H132").Value = 998
$ws.Range("I132").Value = 998
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2994
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -464
